$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.9749262332916433, 1.300361778587099, -6.267426431179062),
    @(-3.844243764877326, 1.833226948976521, -1.409952521324157),
    @(0.4797788858413697, -0.523662269115448, -1.702465817332268),
    @(1.155098915100098, 1.092013478279114, 1.727226853370667),
    @(-1.098365545272828, -0.6193101108074199, 0.1845241859555233),
    @(-0.8518145084381094, -0.03355145454406605, 0.7549576908350003),
    @(0.405293345451355, 0.8384262472391129, 0.3231545425951481),
    @(0.2438197135925255, 0.4860433936119046, -0.09267929568886754),
    @(-0.07322704792022328, 0.1344193816185026, -0.148086081258953),
    @(0.4853796958923352, 0.7483796477317809, -0.1751452423632148)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
